$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 535.2174
$ws.Range("I33").Value = 562.4286
$ws.Range("K33").Value = 562.4286
$ws.Range("M33").Value = -333.4286
$ws.Range("H137").Value = 1895452.2
$ws.Range("I137").Value = 3969648.8
$ws.Range("J137").Value = 1620.6522
$ws.Range("K137").Value = 11908946.4
$ws.Range("L137").Value = 4861.9566
$ws.Range("M137").Value = -11906396.4
$ws.Range("N137").Value = -9961.9566
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 252225
$ws.Range("I122").Value = 1000000
$ws.Range("J122").Value = 2966.6667
$ws.Range("K122").Value = 3000000
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("M122").Value = -2997550
$ws.Range("N122").Value = -13800.0001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 39934
$ws.Range("J4").Value = 39934
$ws.Range("L4").Value = 39934
$ws.Range("N4").Value = -40158
$ws.Range("H22").Value = 430.05554
$ws.Range("I22").Value = 318.75
$ws.Range("J22").Value = 519.1
$ws.Range("K22").Value = 318.75
$ws.Range("L22").Value = 519.1
$ws.Range("M22").Value = 31.25
$ws.Range("N22").Value = -1219.1
$ws.Range("H31").Value = 5535
$ws.Range("I31").Value = 1875.1052
$ws.Range("J31").Value = 7046.696
$ws.Range("K31").Value = 1875.1052
$ws.Range("L31").Value = 7046.696
$ws.Range("M31").Value = -1580.1052
$ws.Range("N31").Value = -7636.696
$ws.Range("H34").Value = 5535
$ws.Range("I34").Value = 1875.1052
$ws.Range("J34").Value = 7046.696
$ws.Range("K34").Value = 1875.1052
$ws.Range("L34").Value = 7046.696
$ws.Range("M34").Value = -1673.1052
$ws.Range("N34").Value = -7450.696
$ws.Range("H59").Value = 16856.572
$ws.Range("I59").Value = 7000
$ws.Range("J59").Value = 18499.334
$ws.Range("K59").Value = 7000
$ws.Range("L59").Value = 18499.334
$ws.Range("M59").Value = -5855
$ws.Range("N59").Value = -20789.334
$ws.Range("H68").Value = 23491.084
$ws.Range("J68").Value = 23491.084
$ws.Range("L68").Value = 23491.084
$ws.Range("N68").Value = -24989.084
$ws.Range("H71").Value = 23491.084
$ws.Range("J71").Value = 23491.084
$ws.Range("L71").Value = 70473.25199999999
$ws.Range("N71").Value = -77961.25199999999
$ws.Range("H74").Value = 20191.4
$ws.Range("J74").Value = 20191.4
$ws.Range("L74").Value = 20191.4
$ws.Range("N74").Value = -21939.4
$ws.Range("H77").Value = 20191.4
$ws.Range("J77").Value = 20191.4
$ws.Range("L77").Value = 60574.2
$ws.Range("N77").Value = -69310.20000000001
$ws.Range("H99").Value = 2422.5112
$ws.Range("J99").Value = 2464.3057
$ws.Range("L99").Value = 2464.3057
$ws.Range("N99").Value = -5460.3057
$ws.Range("H126").Value = 2422.5112
$ws.Range("J126").Value = 2464.3057
$ws.Range("L126").Value = 7392.9171
$ws.Range("N126").Value = -12332.9171
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 471.5238
$ws.Range("I25").Value = 200
$ws.Range("J25").Value = 485.1
$ws.Range("K25").Value = 600
$ws.Range("L25").Value = 1455.3
$ws.Range("M25").Value = -431
$ws.Range("N25").Value = -1793.3
$ws.Range("H30").Value = 471.5238
$ws.Range("I30").Value = 200
$ws.Range("J30").Value = 485.1
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 1455.3
$ws.Range("M30").Value = -498
$ws.Range("N30").Value = -1659.3
$ws.Range("H93").Value = 4973.05
$ws.Range("J93").Value = 4975.737
$ws.Range("L93").Value = 14927.211
$ws.Range("N93").Value = -18671.211
$ws.Range("H140").Value = 1796.7097
$ws.Range("I140").Value = 1187.1578
$ws.Range("J140").Value = 2761.8333
$ws.Range("K140").Value = 3561.4734
$ws.Range("L140").Value = 8285.499899999999
$ws.Range("M140").Value = 1618.5266
$ws.Range("N140").Value = -18645.4999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 12432
$ws.Range("J134").Value = 12432
$ws.Range("L134").Value = 37296
$ws.Range("N134").Value = -42366
$ws.Range("H141").Value = 73286.57000000001
$ws.Range("J141").Value = 73286.57000000001
$ws.Range("L141").Value = 73286.57000000001
$ws.Range("N141").Value = -83646.57000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 8999
$ws.Range("J2").Value = 8999
$ws.Range("L2").Value = 8999
$ws.Range("N2").Value = -9223
$ws.Range("H22").Value = 10570.3
$ws.Range("I22").Value = 750.5
$ws.Range("J22").Value = 13025.25
$ws.Range("K22").Value = 750.5
$ws.Range("L22").Value = 13025.25
$ws.Range("M22").Value = -455.5
$ws.Range("N22").Value = -13615.25
$ws.Range("H27").Value = 10570.3
$ws.Range("I27").Value = 750.5
$ws.Range("J27").Value = 13025.25
$ws.Range("K27").Value = 750.5
$ws.Range("L27").Value = 13025.25
$ws.Range("M27").Value = -643.5
$ws.Range("N27").Value = -13239.25
$ws.Range("H40").Value = 2161.111
$ws.Range("I40").Value = 2208.3333
$ws.Range("J40").Value = 2066.6667
$ws.Range("K40").Value = 2208.3333
$ws.Range("L40").Value = 2066.6667
$ws.Range("M40").Value = -2072.3333
$ws.Range("N40").Value = -2338.6667
$ws.Range("H46").Value = 5444
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5444
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5444
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5820
$ws.Range("H93").Value = 8850.571
$ws.Range("I93").Value = 12233.333
$ws.Range("J93").Value = 2761.6
$ws.Range("K93").Value = 12233.333
$ws.Range("L93").Value = 2761.6
$ws.Range("M93").Value = -10985.333
$ws.Range("N93").Value = -5257.6
$ws.Range("H99").Value = 98285
$ws.Range("J99").Value = 98285
$ws.Range("L99").Value = 98285
$ws.Range("M99").Value = -104275
$ws.Range("H122").Value = 5278.048
$ws.Range("I122").Value = 2300
$ws.Range("J122").Value = 6469.2666
$ws.Range("K122").Value = 6900
$ws.Range("L122").Value = 19407.7998
$ws.Range("M122").Value = -4450
$ws.Range("N122").Value = -24307.7998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3680
$ws.Range("I96").Value = 3009.0908
$ws.Range("K96").Value = 3009.0908
$ws.Range("M96").Value = -1636.0908
$ws.Range("H97").Value = 39653.145
$ws.Range("J97").Value = 39653.145
$ws.Range("L97").Value = 39653.145
$ws.Range("N97").Value = -41635.145
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
